# update: datasets | citation check | document browser
$d = $word.ActiveDocument

# 1) "Following Elsous et al. (Ref-A1B2C3)" -> "Following Elsous et al. (Lee et al., 2020)"
$d.Content.Find.Execute("Elsous et al. (Ref-A1B2C3)", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Elsous et al. (Lee et al., 2020)", 2)

# 2) "patient care (Ref-DJ72KL)" -> "patient care (Ref-u107884)"
$d.Content.Find.Execute("patient care (Ref-DJ72KL)", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "patient care (Ref-u107884)", 2)

# 3) "Health Sciences (Ref-A1B2C3)" -> "Health Sciences (Ref-s219563)"
$d.Content.Find.Execute("Health Sciences (Ref-A1B2C3)", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Health Sciences (Ref-s219563)", 2)

# 4) "successfully (Ref-A1B2C3)" -> "successfully (Pearse et al., 2001)"
$d.Content.Find.Execute("successfully (Ref-A1B2C3)", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "successfully (Pearse et al., 2001)", 2)
